# Weekly update: a new price-report row for the most recent week is
# inserted at row 92 (the top of the "Poroto verde" price history for
# "Feria Lagunitas de Puerto Montt"), pushing the existing rows 92..180
# down to 93..181 and growing the sheet from R180 to R181.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 92; Excel shifts rows 92..180 down to
# 93..181 and carries the column formatting (incl. the date style on D)
# down from the row above, exactly like a manual "Insert Row" in the UI.
$ws.Rows.Item(92).Insert()

# Populate the new row 92 with the new week's record.
$ws.Range("A92").Value = 4
$ws.Range("B92").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C92").Value = "Los Lagos"
$ws.Range("D92").Value = 45233
$ws.Range("E92").Value = 10
$ws.Range("F92").Value = 100112031
$ws.Range("G92").Value = "Poroto verde"
$ws.Range("H92").Value = "Magnum"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 40
$ws.Range("K92").Value = 48000
$ws.Range("L92").Value = 48000
$ws.Range("M92").Value = 48000
$ws.Range("N92").Value = "`$/malla 25 kilos"
$ws.Range("O92").Value = "Perú"
$ws.Range("P92").Value = 1920
$ws.Range("Q92").Value = 25
$ws.Range("R92").Value = "Hortaliza"
